$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.932.73'
$ws.Range('E2').Value = '  +0.37%  '

$ws.Range('D3').Value = '2.362.44'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.57'
$ws.Range('E5').Value = '  +0.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.87'
$ws.Range('E6').Value = '  +0.59%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  -0.46%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.489'
$ws.Range('E9').Value = '  -0.38%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.12'
$ws.Range('E10').Value = '  +0.05%  '

$ws.Range('E11').Value = '  +3.62%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0785'
$ws.Range('E12').Value = '  +0.37%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.32'
$ws.Range('E13').Value = '  -3.20%  '

$ws.Range('E14').Value = '  +0.25%  '

$ws.Range('D15').Value = '2.730.67'
$ws.Range('E15').Value = '  +2.07%  '

$ws.Range('D16').Value = '2.355.74'
$ws.Range('E16').Value = '  +3.52%  '

$ws.Range('E17').Value = '  +0.75%  '

$ws.Range('D18').Value = '42.907.77'
$ws.Range('E18').Value = '  +0.43%  '

$ws.Range('B19').Value = 'InternetComputer(DFINITY)'
$ws.Range('C19').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.88'
$ws.Range('E19').Value = '  -1.64%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.26'
$ws.Range('E20').Value = '  +2.06%  '

$ws.Range('D21').Value = '0.0₃0885'
$ws.Range('E21').Value = '  -0.43%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.93'
$ws.Range('E22').Value = '  +0.35%  '

$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('E24').Value = '  -5.15%  '

$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.43'
$ws.Range('E25').Value = '  +0.79%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.21%  '

$ws.Range('E27').Value = '  +0.74%  '

$ws.Range('E28').Value = '  +0.75%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.32'
$ws.Range('E29').Value = '  +2.24%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.09'
$ws.Range('E30').Value = '  -0.18%  '

$ws.Range('E31').Value = '  -0.08%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.02'
$ws.Range('E32').Value = '  +0.58%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.37'
$ws.Range('E33').Value = '  -1.40%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0717'
$ws.Range('E34').Value = '  +2.69%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '130.79'
$ws.Range('E35').Value = '  -13.13%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.84'
$ws.Range('E36').Value = '  +3.40%  '

$ws.Range('E37').Value = '  +3.49%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.35'
$ws.Range('E38').Value = '  -1.94%  '

$ws.Range('E39').Value = '  -1.88%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.79'
$ws.Range('E40').Value = '  +3.02%  '

$ws.Range('E41').Value = '  -0.59%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.33'
$ws.Range('E42').Value = '  -1.54%  '

$ws.Range('D43').Value = '1.930.44'
$ws.Range('E43').Value = '  +0.72%  '

$ws.Range('E44').Value = '  -0.02%  '

$ws.Range('E45').Value = '  +2.51%  '

$ws.Range('E46').Value = '  -0.72%  '

$ws.Range('E47').Value = '  -9.12%  '

$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.595.48'
$ws.Range('E48').Value = '  +1.99%  '

$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.51'
$ws.Range('E49').Value = '  +2.32%  '

$ws.Range('E50').Value = '  +1.37%  '

$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.49'
$ws.Range('E51').Value = '  -3.21%  '
